# Planeamento de Testes.xlsx - fix mis-placed "X" marks (col C -> col B)
# and tidy the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# --- Row 7: "Adicionar Vencedor" -------------------------------------
# The "X" was typed into C7 instead of B7. Move it to B7 (matching every
# other row in the sheet) and leave C7 present-but-empty, centered.
$ws.Range("B7").Value2 = $ws.Range("C7").Value2
$ws.Range("B7").HorizontalAlignment = $xlCenter
$ws.Range("B7").VerticalAlignment = $xlCenter

$ws.Range("C7").ClearContents()
$ws.Range("C7").HorizontalAlignment = $xlCenter
$ws.Range("C7").VerticalAlignment = $xlCenter
$ws.Range("C7").Font.Name = "Aptos Narrow"

# --- Row 36: same typo ------------------------------------------------
$ws.Range("B36").Value2 = $ws.Range("C36").Value2
$ws.Range("B36").HorizontalAlignment = $xlCenter
$ws.Range("B36").VerticalAlignment = $xlCenter
$ws.Range("C36").Clear()

# --- Row 42: same typo -------------------------------------------------
$ws.Range("B42").Value2 = $ws.Range("C42").Value2
$ws.Range("B42").HorizontalAlignment = $xlCenter
$ws.Range("B42").VerticalAlignment = $xlCenter
$ws.Range("C42").Clear()

# --- Selection / scroll position ---------------------------------------
# Move the cursor to C7 (also drops the stale top-left scroll anchor).
$ws.Range("C7").Select() | Out-Null
